{"js": "// Apply the \"Compact\" paragraph style to the empty paragraphs that sit in\n// the blank \"Due\" cells of the Schedule table (rows for Day 1, 2, 4, 6, 7).\n// These cells currently contain a bare, style-less empty paragraph\n// (<w:p/>); the author gave them the same \"Compact\" style already used by\n// every other populated cell in the table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  return \"no tables found\";\n}\n\n// There is a single table in this document (the Schedule table).\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Gather every cell's first paragraph across the whole table.\nconst cellParagraphs = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  const cell = table.getCell(r, 2); // \"Due\" column\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  cellParagraphs.push({ row: r, paragraphs });\n}\nawait context.sync();\n\n// Load style/text for each paragraph so we can find the untouched,\n// completely empty ones (no explicit style, no text).\ncellParagraphs.forEach(({ paragraphs }) => {\n  paragraphs.items.forEach((p) => p.load(\"style,text\"));\n});\nawait context.sync();\n\nlet changed = 0;\ncellParagraphs.forEach(({ paragraphs }) => {\n  paragraphs.items.forEach((p) => {\n    if (p.style === \"Normal\" && p.text === \"\") {\n      p.style = \"Compact\";\n      changed++;\n    }\n  });\n});\nawait context.sync();\n\nreturn `updated ${changed} empty cell paragraph(s) to Compact style`;\n", "ps1": "# Apply the \"Compact\" paragraph style to the empty paragraphs that sit in\n# the blank \"Due\" cells of the Schedule table (rows for Day 1, 2, 4, 6, 7).\n# These cells currently contain a bare, style-less empty paragraph\n# (<w:p/>); the author gave them the same \"Compact\" style already used by\n# every other populated cell in the table.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$emptyCellMarker = [string][char]13 + [char]7   # Word's cell-end \"end of cell\" marker text\n\n$changed = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $cell = $table.Cell($r, 3)  # \"Due\" column\n    $para = $cell.Range.Paragraphs.Item(1)\n    $isEmpty = ($cell.Range.Text -eq $emptyCellMarker)\n    $styleName = $para.Style.NameLocal\n    if ($isEmpty -and $styleName -eq \"Normal\") {\n        $para.Style = \"Compact\"\n        $changed++\n    }\n}\n\nWrite-Output (\"updated \" + $changed + \" empty cell paragraph(s) to Compact style\")\n"}
